$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.339.40"
$ws.Range("E2").Value = "  -3.41%  "
$ws.Range("D3").Value = "2.599.14"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'568.55"
$ws.Range("E5").Value = "  -4.62%  "
$ws.Range("D6").Value = "'153.09"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.615"
$ws.Range("E8").Value = "  -3.93%  "
$ws.Range("D9").Value = "2.597.47"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("D10").Value = "'0.114"
$ws.Range("E10").Value = "  -8.38%  "
$ws.Range("D11").Value = "'5.77"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.156"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "'0.375"
$ws.Range("E13").Value = "  -5.41%  "
$ws.Range("D14").Value = "'27.77"
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("D15").Value = "3.070.81"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = "  -8.20%  "
$ws.Range("D17").Value = "63.282.00"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").Value = "2.595.65"
$ws.Range("E18").Value = "  -3.33%  "
$ws.Range("D19").Value = "'11.85"
$ws.Range("E19").Value = "  -4.80%  "
$ws.Range("D20").Value = "'7.42"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "'4.45"
$ws.Range("E21").Value = "  -6.38%  "
$ws.Range("D22").Value = "'338.25"
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'66.96"
$ws.Range("E24").Value = "  -3.54%  "
$ws.Range("D25").Value = "'1.78"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("D26").Value = "'0.0000104"
$ws.Range("E26").Value = "  -7.19%  "
$ws.Range("D27").Value = "'9.00"
$ws.Range("E27").Value = "  -6.34%  "
$ws.Range("D28").Value = "'574.80"
$ws.Range("E28").Value = "  +3.32%  "
$ws.Range("D29").Value = "'1.53"
$ws.Range("E29").Value = "  -4.89%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "'0.159"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").Value = "'7.72"
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("E33").Value = "  -4.81%  "
$ws.Range("D34").Value = "'1.70"
$ws.Range("E34").Value = "  -5.44%  "
$ws.Range("D35").Value = "'6.46"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("D36").Value = "'5.31"
$ws.Range("E36").Value = "  -2.47%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.397"
$ws.Range("E38").Value = "  -5.51%  "
$ws.Range("D39").Value = "'19.51"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").Value = "'154.48"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "'1.84"
$ws.Range("E41").Value = "  -5.72%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'41.51"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").Value = "'2.45"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'156.52"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("D46").Value = "'23.13"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'3.81"
$ws.Range("E47").Value = "  -6.35%  "
$ws.Range("D48").Value = "'0.0576"
$ws.Range("E48").Value = "  -6.02%  "
$ws.Range("D49").Value = "'0.624"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("D50").Value = "'0.0989"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").Value = "'0.0243"
$ws.Range("E51").Value = "  -5.45%  "
